$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '07/14/2022'
$ws.Range("J2").Value = 65
$ws.Range("K2").Value = 14000
$ws.Range("L2").Value = 14000
$ws.Range("M2").Value = 14000
$ws.Range("N2").Value = '$/caja 15 kilos'
$ws.Range("P2").Value = 933
$ws.Range("Q2").Value = 15

$ws.Range("D3").Value = '07/14/2022'
$ws.Range("J3").Value = 68
$ws.Range("K3").Value = 12000
$ws.Range("L3").Value = 12000
$ws.Range("M3").Value = 12000
$ws.Range("N3").Value = '$/caja 15 kilos'
$ws.Range("P3").Value = 800
$ws.Range("Q3").Value = 15

$ws.Range("D4").Value = '03/07/2023'
$ws.Range("J4").Value = 56
$ws.Range("K4").Value = 13000
$ws.Range("L4").Value = 13000
$ws.Range("M4").Value = 13000
$ws.Range("N4").Value = '$/bandeja 18 kilos'
$ws.Range("P4").Value = 722
$ws.Range("Q4").Value = 18

$ws.Range("D5").Value = '02/15/2021'
$ws.Range("I5").Value = 'Primera'
$ws.Range("J5").Value = 60
$ws.Range("K5").Value = 13000
$ws.Range("L5").Value = 13000
$ws.Range("M5").Value = 13000
$ws.Range("N5").Value = '$/bandeja 18 kilos'
$ws.Range("P5").Value = 722
$ws.Range("Q5").Value = 18

$ws.Range("D6").Value = '02/15/2021'
$ws.Range("I6").Value = 'Segunda'
$ws.Range("J6").Value = 50
$ws.Range("K6").Value = 10000
$ws.Range("L6").Value = 10000
$ws.Range("M6").Value = 10000
$ws.Range("P6").Value = 556

$ws.Range("I7").Value = 'Primera'
$ws.Range("J7").Value = 87
$ws.Range("K7").Value = 22000
$ws.Range("L7").Value = 22000
$ws.Range("M7").Value = 22000
$ws.Range("P7").Value = 1222

$ws.Range("D8").Value = '12/06/2021'
$ws.Range("I8").Value = 'Segunda'
$ws.Range("K8").Value = 20000
$ws.Range("L8").Value = 20000
$ws.Range("M8").Value = 20000
$ws.Range("P8").Value = 1111

$ws.Range("D9").Value = '03/06/2023'
$ws.Range("I9").Value = 'Primera'
$ws.Range("J9").Value = 75
$ws.Range("K9").Value = 13000
$ws.Range("L9").Value = 13000
$ws.Range("M9").Value = 13000
$ws.Range("P9").Value = 722

$ws.Range("D10").Value = '03/06/2023'
$ws.Range("I10").Value = 'Segunda'
$ws.Range("J10").Value = 56
$ws.Range("K10").Value = 9000
$ws.Range("L10").Value = 9000
$ws.Range("M10").Value = 9000
$ws.Range("P10").Value = 500

$ws.Range("D11").Value = '02/08/2021'
$ws.Range("J11").Value = 80
$ws.Range("K11").Value = 14000
$ws.Range("L11").Value = 14000
$ws.Range("M11").Value = 14000
$ws.Range("P11").Value = 778

$ws.Range("D12").Value = '02/08/2021'
$ws.Range("J12").Value = 70
$ws.Range("K12").Value = 12000
$ws.Range("L12").Value = 12000
$ws.Range("M12").Value = 12000
$ws.Range("P12").Value = 667

$ws.Range("D13").Value = '02/08/2021'
$ws.Range("I13").Value = 'Tercera'
$ws.Range("J13").Value = 60
$ws.Range("K13").Value = 10000
$ws.Range("L13").Value = 10000
$ws.Range("M13").Value = 10000
$ws.Range("N13").Value = '$/bandeja 18 kilos'
$ws.Range("P13").Value = 556
$ws.Range("Q13").Value = 18

$ws.Range("D14").Value = '02/11/2021'
$ws.Range("I14").Value = 'Primera'
$ws.Range("J14").Value = 90
$ws.Range("K14").Value = 13000
$ws.Range("L14").Value = 13000
$ws.Range("M14").Value = 13000
$ws.Range("N14").Value = '$/bandeja 18 kilos'
$ws.Range("P14").Value = 722
$ws.Range("Q14").Value = 18

$ws.Range("D15").Value = '02/11/2021'
$ws.Range("I15").Value = 'Segunda'
$ws.Range("J15").Value = 80
$ws.Range("K15").Value = 11000
$ws.Range("L15").Value = 11000
$ws.Range("M15").Value = 11000
$ws.Range("P15").Value = 611

$ws.Range("D16").Value = '08/16/2021'
$ws.Range("I16").Value = 'Primera'
$ws.Range("J16").Value = 75
$ws.Range("K16").Value = 18000
$ws.Range("L16").Value = 18000
$ws.Range("M16").Value = 18000
$ws.Range("N16").Value = '$/caja 15 kilos'
$ws.Range("P16").Value = 1200
$ws.Range("Q16").Value = 15

$ws.Range("D17").Value = '08/16/2021'
$ws.Range("I17").Value = 'Segunda'
$ws.Range("J17").Value = 50
$ws.Range("K17").Value = 12000
$ws.Range("L17").Value = 12000
$ws.Range("M17").Value = 12000
$ws.Range("N17").Value = '$/caja 15 kilos'
$ws.Range("P17").Value = 800
$ws.Range("Q17").Value = 15
